$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.624.79"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "3.009.79"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.86"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.68"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.54"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("E11").Value = "  +4.22%  "
$ws.Range("D12").Value = "3.524.65"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("E14").Value = "  +4.09%  "
$ws.Range("E15").Value = "  +7.28%  "
$ws.Range("D16").Value = "57.567.81"
$ws.Range("E16").Value = "  +1.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.20"
$ws.Range("E17").Value = "  +6.02%  "
$ws.Range("D18").Value = "3.008.07"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.84"
$ws.Range("E19").Value = "  +3.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.98"
$ws.Range("E20").Value = "  +2.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.23"
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E24").Value = "  +3.51%  "
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("E27").Value = "  +2.27%  "
$ws.Range("E28").Value = "  +3.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("E29").Value = "  +4.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.81"
$ws.Range("E30").Value = "  +2.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  -5.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.59"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.72"
$ws.Range("E33").Value = "  +5.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "154.34"
$ws.Range("E34").Value = "  -0.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.86"
$ws.Range("E35").Value = "  +4.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.28"
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "24.42"
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "3.042.41"
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.34"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.84"
$ws.Range("E41").Value = "  +6.93%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "2.244.08"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.986"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.02"
$ws.Range("E47").Value = "  +4.93%  "
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.37"
$ws.Range("E49").Value = "  +2.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.86"
$ws.Range("E50").Value = "  -6.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0893"
$ws.Range("E51").Value = "  +2.73%  "
